$d = $word.ActiveDocument

# Locate the "Improvements" heading paragraph (start of the section being
# reverted/removed) and the paragraph holding the last bullet ("Make High
# poly level feel colder..."), which also carries the trailing _GoBack
# bookmark right before its own paragraph mark.
$startIndex = -1
$endIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.Trim()
    if ($startIndex -eq -1 -and $t -eq "Improvements") {
        $startIndex = $i
    }
    if ($t -like "Make High poly level feel colder*") {
        $endIndex = $i
    }
}

if ($startIndex -ne -1 -and $endIndex -ne -1 -and $endIndex -gt $startIndex) {
    # Step 1: remove the "Improvements" heading through the second-to-last
    # bullet ("Make low poly level night time..."), i.e. everything up to
    # (but not including) the paragraph that holds the trailing bookmark.
    $beforeLastPara = $d.Paragraphs.Item($endIndex - 1)
    $startPara = $d.Paragraphs.Item($startIndex)
    $rng = $d.Range($startPara.Range.Start, $beforeLastPara.Range.End)
    $rng.Delete()

    # Step 2: remove the visible text of the remaining bullet paragraph (the
    # one that used to read "Make High poly level feel colder, with more
    # snow than the other level. "), leaving only its paragraph mark plus
    # the bookmark that sits right before it.
    $lastPara = $d.Paragraphs.Item($startIndex)
    $textOnly = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
    $textOnly.Delete()

    # Step 3: merge that now-empty paragraph back into the previous one (the
    # "UI bloops for selecting menu options." paragraph) by deleting the
    # preceding paragraph's own paragraph mark. This keeps the bookmark
    # intact and moves it to sit right after " for selecting menu options. "
    # in the same paragraph, matching the reverted document structure.
    $prevPara = $d.Paragraphs.Item($startIndex - 1)
    $markRng = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End)
    $markRng.Delete()
}
